# Insert two new data rows at row 316 (this shifts the existing rows 316-361
# down to 318-363, extending the used range from A1:R361 to A1:R363).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(316).Insert()
$ws.Rows.Item(316).Insert()

# New row 316
$ws.Range("A316").Value = 8
$ws.Range("B316").Value = "Terminal La Palmera de La Serena"
$ws.Range("C316").Value = "Coquimbo"
$ws.Range("D316").Value = 44951
$ws.Range("D316").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E316").Value = 4
$ws.Range("F316").Value = 100112021
$ws.Range("G316").Value = "Ají"
$ws.Range("H316").Value = "Inferno"
$ws.Range("I316").Value = "Primera"
$ws.Range("J316").Value = 460
$ws.Range("K316").Value = 13500
$ws.Range("L316").Value = 14000
$ws.Range("M316").Value = 13750
$ws.Range("N316").Value = "$/caja 15 kilos"
$ws.Range("O316").Value = "Provincia de Limarí"
$ws.Range("P316").Value = 917
$ws.Range("Q316").Value = 15
$ws.Range("R316").Value = "Hortaliza"

# New row 317
$ws.Range("A317").Value = 8
$ws.Range("B317").Value = "Terminal La Palmera de La Serena"
$ws.Range("C317").Value = "Coquimbo"
$ws.Range("D317").Value = 44951
$ws.Range("D317").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E317").Value = 4
$ws.Range("F317").Value = 100112021
$ws.Range("G317").Value = "Ají"
$ws.Range("H317").Value = "Inferno"
$ws.Range("I317").Value = "Segunda"
$ws.Range("J317").Value = 360
$ws.Range("K317").Value = 9500
$ws.Range("L317").Value = 10000
$ws.Range("M317").Value = 9750
$ws.Range("N317").Value = "$/caja 15 kilos"
$ws.Range("O317").Value = "Provincia de Limarí"
$ws.Range("P317").Value = 650
$ws.Range("Q317").Value = 15
$ws.Range("R317").Value = "Hortaliza"
